$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-24 Friday" "2025-01-25 Saturday"

Replace-Text "995÷3=" "721÷6="
Replace-Text "851÷4=" "545÷6="
Replace-Text "356÷2=" "719÷9="
Replace-Text "847÷8=" "911÷5="
Replace-Text "279÷9=" "616÷5="

Replace-Text "953÷5=" "951÷4="
Replace-Text "311÷8=" "741÷7="
Replace-Text "942÷2=" "896÷6="
Replace-Text "400÷5=" "218÷3="
Replace-Text "966÷8=" "875÷7="

Replace-Text "316÷7=" "800÷3="
Replace-Text "392÷9=" "937÷3="
Replace-Text "194÷6=" "706÷4="
Replace-Text "927÷5=" "917÷8="
Replace-Text "584÷3=" "450÷4="

Replace-Text "348÷5=" "868÷4="
Replace-Text "392÷2=" "285÷6="
Replace-Text "614÷2=" "697÷9="
Replace-Text "211÷6=" "308÷5="
Replace-Text "126÷4=" "253÷9="

Replace-Text "991÷3=" "892÷6="
Replace-Text "303÷9=" "943÷5="
Replace-Text "505÷6=" "920÷5="
Replace-Text "537÷5=" "513÷7="
Replace-Text "778÷8=" "715÷3="
